$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data sheet")
$ws.Range("G2").Value = 2433
$ws.Range("A3:L3").Copy()
$ws.Range("A4").PasteSpecial()
$ws.Range("K4").Value = "R2, R3, r4, r5"
$ws.Range("K4").Select()
